# Insert a new weekly price record at row 502 ("Fruta / hortaliza, semanal").
# Excel shifts every following row (502-567) down by one (503-568), which is
# exactly what the target diff shows (each old row's data reappears one row
# lower, and the former last row 567 is duplicated into the new last row 568).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 502:567 down to 503:568, duplicating formatting (incl. the date
# style on column D) from the row being pushed down, same as a manual
# "Insert" in Excel.
$ws.Rows.Item(502).Insert()

# Populate the newly inserted row 502 with the new record's data.
$ws.Cells.Item(502, 1).Value  = 8
$ws.Cells.Item(502, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(502, 3).Value  = "Coquimbo"
$ws.Cells.Item(502, 4).Value  = 45154
$ws.Cells.Item(502, 5).Value  = 4
$ws.Cells.Item(502, 6).Value  = 100112032
$ws.Cells.Item(502, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(502, 8).Value  = "Sin especificar"
$ws.Cells.Item(502, 9).Value  = "Primera"
$ws.Cells.Item(502, 10).Value = 440
$ws.Cells.Item(502, 11).Value = 14000
$ws.Cells.Item(502, 12).Value = 15000
$ws.Cells.Item(502, 13).Value = 14500
$ws.Cells.Item(502, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(502, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(502, 16).Value = 290
$ws.Cells.Item(502, 17).Value = 50
$ws.Cells.Item(502, 18).Value = "Hortaliza"
